{"js": "const body = context.document.body;\n\n// Locate the paragraph that currently holds the M2Doc field\n// \"m:'Mona_Lisa.jpg'.asImage().fit(300, 150, false)\" (a real Word field:\n// fldChar begin / instrText runs / fldChar end).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].fields.load(\"items/code\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const fields = paragraphs.items[i].fields;\n  for (let j = 0; j < fields.items.length; j++) {\n    const code = fields.items[j].code || \"\";\n    if (code.indexOf(\"asImage\") !== -1) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n  if (target) {\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the expected M2Doc image field.\");\n}\n\n// Replace the whole paragraph (field delimiters included) with plain text\n// runs: literal \"{\" / \"}\" instead of field begin/end markers, and the same\n// orange accent-6 color kept on the templating expression itself.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t>{</w:t></w:r>\n<w:r><w:t>m</w:t></w:r>\n<w:r><w:t>:</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>'</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>Mona_Lisa.jpg</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>'.asImage()</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.fit(</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>300</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t xml:space=\"preserve\">, </w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>150</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>, false</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>)</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the M2Doc field \"m:'Mona_Lisa.jpg'.asImage().fit(300, 150, false)\"\n# which is currently stored as a real Word field (fldChar begin/instrText/fldChar end).\n$target = $null\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n    $candidate = $d.Fields.Item($i)\n    if ($candidate.Code.Text -match \"asImage\") {\n        $target = $candidate\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the expected M2Doc image field.\"\n}\n\n# Find the index of the paragraph that hosts the field (scanning avoids a\n# collapsed-range Paragraphs.Item(1) quirk that resolves to the wrong paragraph).\n$fieldStart = $target.Code.Start\n$paraIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $pr = $d.Paragraphs.Item($i).Range\n    if ($fieldStart -ge $pr.Start -and $fieldStart -le $pr.End) {\n        $paraIndex = $i\n        break\n    }\n}\n\n# Remove the field (fldChar begin / instrText runs / fldChar end) entirely.\n$target.Delete()\n\n# Re-acquire the (now empty) paragraph range, keeping the trailing paragraph mark out.\n$para = $d.Paragraphs.Item($paraIndex)\n$r = $para.Range\n$r.End = $r.End - 1\n\n# Rebuild the same text as plain runs (curly braces instead of field delimiters),\n# keeping the orange accent-6 color on the templating expression itself.\n$xmlFragment = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t>{</w:t></w:r>\n<w:r><w:t>m</w:t></w:r>\n<w:r><w:t>:</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>'</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>Mona_Lisa.jpg</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>'.asImage()</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.fit(</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>300</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t xml:space=\"preserve\">, </w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>150</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>, false</w:t></w:r>\n<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>)</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$r.InsertXML($xmlFragment)\n"}
